$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new values in column I, mirroring column E's style (centered, style index 2)
$ws.Range("I3").Value = 6866
$ws.Range("I5").Value = 270
$ws.Range("I6").Value = 60
$ws.Range("I12").Value = 910

# Apply the same formatting as column E (style index "2" -> horizontal center)
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("I5").HorizontalAlignment = -4108
$ws.Range("I6").HorizontalAlignment = -4108
$ws.Range("I12").HorizontalAlignment = -4108

# Update selection
$ws.Range("C4").Select()
